$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.217.86'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.643.96'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.23'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.95'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.873.83'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.655.69'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.545'
$ws.Range("E15").Value = '  +3.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.43'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.199.25'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("E21").Value = '  +3.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.55'
$ws.Range("E22").Value = '  +4.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.40'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +1.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.55'
$ws.Range("E26").Value = '  +2.58%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.262.09'
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("E42").Value = '  +6.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.34'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.783.87'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.83'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.66'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0514'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  +1.90%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0974'
$ws.Range("E51").Value = '  +0.26%  '
